# Auto-generated Word COM-interop script
# Applies targeted paragraph rewrites (proofErr spell-check markers + text fixes)
$d = $word.ActiveDocument

# Paragraph 23 (paraId 1A8F6366)
$xml23 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1A8F6366" w14:textId="68941124" w:rsidR="00FC3CB6" w:rsidRDefault="00FC3CB6" w:rsidP="006C09EA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Data script from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> backup?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(23).Range.InsertXML($xml23)

# Paragraph 22 (paraId 182143DD)
$xml22 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="182143DD" w14:textId="2E8F9A86" w:rsidR="00DF3273" w:rsidRDefault="004241E1" w:rsidP="006C09EA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Smallie: ability to give an AWS secret name in command line , instead of full </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>connstr</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(22).Range.InsertXML($xml22)

# Paragraph 20 (paraId 2B80ACCA)
$xml20 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2B80ACCA" w14:textId="54CD1BF9" w:rsidR="002D77F1" w:rsidRPr="00E628B0" w:rsidRDefault="002D77F1" w:rsidP="002D77F1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr><w:rPr><w:color w:val="80340D" w:themeColor="accent2" w:themeShade="80"/></w:rPr></w:pPr><w:r w:rsidRPr="00E628B0"><w:rPr><w:color w:val="80340D" w:themeColor="accent2" w:themeShade="80"/></w:rPr><w:t xml:space="preserve">Biggest by far: when running, also generate a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00E628B0"><w:rPr><w:color w:val="80340D" w:themeColor="accent2" w:themeShade="80"/></w:rPr><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00E628B0"><w:rPr><w:color w:val="80340D" w:themeColor="accent2" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> of the database diffs. Then Electron GUI showing grid</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(20).Range.InsertXML($xml20)

# Paragraph 15 (paraId 7DA8A916)
$xml15 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7DA8A916" w14:textId="565FE805" w:rsidR="000A39FE" w:rsidRDefault="000A39FE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Smallie: title of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>server.dbname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and timestamp ( as usual, will be in comments or printed out)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(15).Range.InsertXML($xml15)

# Paragraph 12 (paraId 0FFE1765)
$xml12 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0FFE1765" w14:textId="77777777" w:rsidR="00D66F8E" w:rsidRDefault="00D66F8E" w:rsidP="00D66F8E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t>“</w:t></w:r><w:r w:rsidRPr="005A0BBA"><w:t xml:space="preserve">bubble up </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005A0BBA"><w:t>differenecs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005A0BBA"><w:t xml:space="preserve"> to table</w:t></w:r><w:r><w:t>”: why don’t I have it in my new script (only on old)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(12).Range.InsertXML($xml12)

# Paragraph 11 (paraId 37C80895)
$xml11 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="37C80895" w14:textId="77777777" w:rsidR="00D66F8E" w:rsidRDefault="00D66F8E" w:rsidP="00D66F8E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r w:rsidRPr="00CF1C63"><w:t xml:space="preserve">#! title not clear. why "the code" is different from "overall code"? do full alignment of everything, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CF1C63"><w:t>wheres</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CF1C63"><w:t xml:space="preserve"> the END of this one? what does this block achieves? </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(11).Range.InsertXML($xml11)

# Paragraph 9 (paraId 0A5B4DA3)
$xml9 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0A5B4DA3" w14:textId="72D62FBB" w:rsidR="00744209" w:rsidRDefault="00744209" w:rsidP="00744209"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t>See where its searching for ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00744209"><w:t>data_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">’ in code. Needs to be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_type_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. make sure </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>its</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> not a bug</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(9).Range.InsertXML($xml9)

# Paragraph 8 (paraId 468ED549)
$xml8 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="468ED549" w14:textId="7B508035" w:rsidR="00B10205" w:rsidRDefault="00B10205" w:rsidP="007C0FCA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">How do we encrypt password in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r><w:r w:rsidR="00941D93"><w:t xml:space="preserve"> Consult </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00941D93"><w:t>claude</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00941D93"><w:t xml:space="preserve"> for ideas. Maybe command line override: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00941D93"><w:t>pwd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00941D93"><w:t xml:space="preserve">, or entire </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00941D93"><w:t>connstr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00941D93"><w:t xml:space="preserve">. Ask </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00941D93"><w:t>claude</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00941D93"><w:t xml:space="preserve"> about that option as well</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(8).Range.InsertXML($xml8)

# Paragraph 7 (paraId 040DA41B)
$xml7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="040DA41B" w14:textId="7E22B6CB" w:rsidR="00A542AB" w:rsidRPr="001301D1" w:rsidRDefault="00A542AB" w:rsidP="001301D1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="10"/></w:numPr><w:rPr><w:color w:val="EE0000"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="EE0000"/></w:rPr><w:t>Then:</w:t></w:r><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve"> finish now the mechanism of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve"> to say which </w:t></w:r><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve">specific </w:t></w:r><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve">schema </w:t></w:r><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve">and data </w:t></w:r><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t>tables to script</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(7).Range.InsertXML($xml7)

# Paragraph 6 (paraId 627E362B)
$xml6 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="627E362B" w14:textId="042D2FC8" w:rsidR="001301D1" w:rsidRDefault="001301D1" w:rsidP="001301D1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="10"/></w:numPr><w:rPr><w:color w:val="EE0000"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="EE0000"/></w:rPr><w:t>RRN:</w:t></w:r><w:r w:rsidRPr="001301D1"><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A542AB"><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t xml:space="preserve">loading data tables! </w:t></w:r><w:r><w:rPr><w:color w:val="EE0000"/></w:rPr><w:t>Load them all but final script has errors</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(6).Range.InsertXML($xml6)

# Paragraph 4 (paraId 03CE85C9)
$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="03CE85C9" w14:textId="0ECE0602" w:rsidR="00901156" w:rsidRDefault="00901156" w:rsidP="00901156"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Now</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r w:rsidRPr="00901156"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">The full load from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>connstr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> given in command line</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(4).Range.InsertXML($xml4)

# Paragraph 2 (paraId 3B6DBA27)
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3B6DBA27" w14:textId="345F6516" w:rsidR="00DE39F8" w:rsidRPr="00FA402B" w:rsidRDefault="00DE39F8" w:rsidP="0017308A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r w:rsidRPr="00FA402B"><w:t xml:space="preserve">when not in full DML mode, but there should still be some text about data </w:t></w:r><w:r w:rsidR="00BF1B0A" w:rsidRPr="00FA402B"><w:t>that was just changed (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00BF1B0A" w:rsidRPr="00FA402B"><w:t>recordcount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00BF1B0A" w:rsidRPr="00FA402B"><w:t xml:space="preserve"> or some)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(2).Range.InsertXML($xml2)

# Paragraph 1 (paraId 398F5278)
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="398F5278" w14:textId="56764519" w:rsidR="005E73B4" w:rsidRDefault="001E7502" w:rsidP="001E7502"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>ContextFreeSQL</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(1).Range.InsertXML($xml1)

Write-Output "done"